# Add a new recipe row (row 5: "PotatoCurry") to the sheet, matching the
# shape/formatting of the existing data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of data (row 5) ---
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "PotatoCurry"
$ws.Range("C5").Value = "tomatoes-2,onion-2,greenchills-3,drychiiles-3,mustardseeds,coriander,cuminseeds,uraddaal"
$ws.Range("D5").Value = "tomatoes-2,onion-2,greenchills-3,drychiiles-3,mustardseeds,coriander,cuminseeds,uraddaal"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = $true
$ws.Range("I5").Value = $false
$ws.Range("L5").Value = 4

# Ingredients/Method columns use the wrap-text cell style (same as row 2's
# C/D cells) throughout the sheet.
$ws.Range("C5").WrapText = $true
$ws.Range("D5").WrapText = $true

# The new row is tall (wrapped long text), like the other data rows.
$ws.Rows.Item(5).RowHeight = 90

# Column L (ReciepeID) gets an explicit width once a value appears in it.
$ws.Columns.Item(12).ColumnWidth = 15.68

# Move the selection to the new row and scroll it into view.
$ws.Range("A5").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
